$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data table (rows 16-23), columns: C=Doc, D=Nombre, E=Periodo, F=ValorMora, G=SalarioBasico
$data = @(
    @{Row=16; C="45537049";   D="MARIA IRENE SEGRERA FUENMAYOR"; E="1812"; F=40000; G=1000000},
    @{Row=17; C="91529249";   D="MARLON VICENTE BACCA MEDINA";   E="1812"; F=60000; G=1500000},
    @{Row=18; C="45560342";   D="VANESSA ROJAS OLMOS";           E="1812"; F=60000; G=0},
    @{Row=19; C="45560342";   D="VANESSA ROJAS OLMOS";           E="1810"; F=60000; G=0},
    @{Row=20; C="1047410598"; D="WILFRIDO MEDINA CEBALLOS";      E="1812"; F=33125; G=800000},
    @{Row=21; C="1140830535"; D="PEDRO JAVIER LLANOS MORALES";   E="1812"; F=48000; G=1200000},
    @{Row=22; C="1047451646"; D="RUBEN DARIO GOMEZ HERRERA";     E="1812"; F=31249; G=800000},
    @{Row=23; C="3809345";    D="TUBAL PADILLA SIMANCAS";        E="1812"; F=40000; G=1000000}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

# Columns B to J get auto-fit to reflect the new bestFit widths after content changes
$ws.Range("B:J").Columns.AutoFit()
